$d = $word.ActiveDocument

# Map the old numeric placeholder tokens ({{idN}}) to the new, descriptive
# merge-field names used by the updated template.
$replacements = @(
    @("{{id1}}",  "{{ID}}"),
    @("{{id2}}",  "{{FullName}}"),
    @("{{id3}}",  "{{StudNum}}"),
    @("{{id4}}",  "{{Theme}}"),
    @("{{id5}}",  "{{SuData}}"),
    @("{{id6}}",  "{{SuName}}"),
    @("{{id11}}", "{{Questioner1}}"),
    @("{{id12}}", "{{Question1}}"),
    @("{{id13}}", "{{Questioner2}}"),
    @("{{id14}}", "{{Question2}}"),
    @("{{id15}}", "{{Questioner3}}"),
    @("{{id16}}", "{{Question3}}"),
    @("{{id18}}", "{{fullName}}"),
    @("{{id17}}", "{{Score}}"),
    @("{{id23}}", "{{Language}}"),
    @("{{id19}}", "{{IndividualOpinion}}")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
